$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-prefix the existing text labels with "%"
$ws.Range("A1").Value = "%foo"
$ws.Range("B1").Value = "%bar"

# Convert the numeric values to "%"-prefixed text (now stored as shared
# strings rather than numbers)
$ws.Range("A2").Value = "%17"
$ws.Range("B2").Value = "%42"

# The four cells now carry their own explicit cell format (distinct from
# the sheet default style) - align this with how the XLSX utility layer
# stamps written cells: explicit horizontal alignment + explicit
# protection/locked state.
$r = $ws.Range("A1:B2")
$r.HorizontalAlignment = -4131
$r.Locked = $true

# Move the active selection from A3 to B3
$ws.Range("B3").Select()
